$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 106: Resistencia vs Tacuary -------------------------------------
$ws.Cells.Item(106, 1).Value  = 105
$ws.Cells.Item(106, 2).Value  = "paraguay"
$ws.Cells.Item(106, 3).Value  = "primera-division"
$ws.Cells.Item(106, 4).NumberFormat = "@"
$ws.Cells.Item(106, 4).Value  = "2023"
$ws.Cells.Item(106, 5).Value  = 45232.91666666666
$ws.Cells.Item(106, 6).Value  = "Resistencia"
$ws.Cells.Item(106, 7).Value  = 1
$ws.Cells.Item(106, 8).Value  = "Tacuary"
$ws.Cells.Item(106, 9).Value  = 3
$ws.Cells.Item(106, 10).Value = 2.1
$ws.Cells.Item(106, 11).Value = "31/10/2023 00:42"
$ws.Cells.Item(106, 12).Value = 2.53
$ws.Cells.Item(106, 13).Value = "02/11/2023 21:58"
$ws.Cells.Item(106, 14).Value = 3.5
$ws.Cells.Item(106, 15).Value = "31/10/2023 00:42"
$ws.Cells.Item(106, 16).Value = 3.46
$ws.Cells.Item(106, 17).Value = "02/11/2023 21:57"
$ws.Cells.Item(106, 18).Value = 3.56
$ws.Cells.Item(106, 19).Value = "31/10/2023 00:42"
$ws.Cells.Item(106, 20).Value = 2.86
$ws.Cells.Item(106, 21).Value = "02/11/2023 21:58"
$ws.Cells.Item(106, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/resistencia-tacuary/EBMT9J5K/"

# --- Row 107: Sportivo Trinidense vs Libertad Asuncion -------------------
$ws.Cells.Item(107, 1).Value  = 106
$ws.Cells.Item(107, 2).Value  = "paraguay"
$ws.Cells.Item(107, 3).Value  = "primera-division"
$ws.Cells.Item(107, 4).NumberFormat = "@"
$ws.Cells.Item(107, 4).Value  = "2023"
$ws.Cells.Item(107, 5).Value  = 45233.02083333334
$ws.Cells.Item(107, 6).Value  = "Sportivo Trinidense"
$ws.Cells.Item(107, 7).Value  = 1
$ws.Cells.Item(107, 8).Value  = "Libertad Asuncion"
$ws.Cells.Item(107, 9).Value  = 1
$ws.Cells.Item(107, 10).Value = 4.31
$ws.Cells.Item(107, 11).Value = "31/10/2023 00:42"
$ws.Cells.Item(107, 12).Value = 4.64
$ws.Cells.Item(107, 13).Value = "03/11/2023 00:24"
$ws.Cells.Item(107, 14).Value = 3.78
$ws.Cells.Item(107, 15).Value = "31/10/2023 00:42"
$ws.Cells.Item(107, 16).Value = 3.86
$ws.Cells.Item(107, 17).Value = "03/11/2023 00:27"
$ws.Cells.Item(107, 18).Value = 1.82
$ws.Cells.Item(107, 19).Value = "31/10/2023 00:42"
$ws.Cells.Item(107, 20).Value = 1.78
$ws.Cells.Item(107, 21).Value = "03/11/2023 00:21"
$ws.Cells.Item(107, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/sportivo-trinidense-libertad-asuncion/roRX8wLQ/"

# --- Replicate the formatting of the previous data row (row 105) onto the
#     two new rows: column A keeps the bordered/bold "index" style and
#     column E keeps the datetime display style; everything else reverts
#     to the default (General) style, exactly like every other data row.
$ws.Range("A105:V105").Copy() | Out-Null
$ws.Range("A106:V107").PasteSpecial(-4122) | Out-Null
